$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells - copy style from an existing header cell (H1) so the new
# headers match the bold/centered/bordered look of the other headers.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for columns I (I0) and J (IF), rows 2-13
$dataI = @(5, 5, 7, 3, 2, 8, 4, 6, 4, 3, 2, 1)
$dataJ = @(6, 6, 7, 5, 5, 8, 6, 9, 7, 4, 3, 2)

for ($i = 0; $i -lt $dataI.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $dataI[$i]
    $ws.Cells.Item($row, 10).Value = $dataJ[$i]
}
